# Update the "想去人数" (interested-people count) figures that were refreshed
# by the scraper run (gh-pages output generated at 456a3b4).
#
# Sheet "展览" (Exhibitions):
#   F2: 408  -> 411
#   F3: 5114 -> 5137
#   F4: 44   -> 46
#   F5: 46   -> 48
#
# Sheet "全部类型" (All types):
#   F2: 408  -> 411
#   F3: 5114 -> 5137
#   F5: 44   -> 46
#   F6: 46   -> 48

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 411
$wsExhibit.Range("F3").Value = 5137
$wsExhibit.Range("F4").Value = 46
$wsExhibit.Range("F5").Value = 48

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 411
$wsAll.Range("F3").Value = 5137
$wsAll.Range("F5").Value = 46
$wsAll.Range("F6").Value = 48
